# Edit the document's primary header (word/header2.xml):
#   - merge the "Don't" / " Panic Room" runs into a single run, drop the
#     manual line break that used to separate the title from the logo/tagline,
#   - replace the struck-through separator "?\u00a0\u00a0KI" (proofErr-wrapped,
#     gramStart/gramEnd) with a plain " \u2013 KI" built from clean runs,
#   - tag every (re)written run/paragraph with English (UK) proofing language,
#   - append a second header paragraph with the "02Station Prompt-Anfaenge
#     Elterngespraeche" title, wrapping the German words "Anfaenge" and
#     "Elterngespraeche" in spell-check proofErr markers, just as Word itself
#     would emit after a live edit + spell pass.
# This mirrors exactly what a user would get by retyping that header text in
# Word, so we drive it through Range.InsertXML with a literal WordprocessingML
# fragment (the precise run/proofErr layout isn't reachable purely through
# Find/Replace, since Word collapses/re-splits runs on its own).

$d = $word.ActiveDocument

# Locate the header that holds the "Don't Panic Room" title - this is the
# document's default/primary header (wdHeaderFooterIndex = 1).
$targetHeader = $null
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $h = $d.Sections.Item($s).Headers.Item(1)
    if ($h.Exists -and $h.Range.Text -like "*Panic Room*") {
        $targetHeader = $h
        break
    }
}
if ($targetHeader -eq $null) {
    $targetHeader = $d.Sections.Item(1).Headers.Item(1)
}

$rng = $targetHeader.Range.Duplicate

# Shared run/paragraph properties: Calibri Light 8pt, tagged English (UK).
$rPr = '<w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/>' + `
       '<w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/>'
$pPr = '<w:pPr><w:pStyle w:val="Kopfzeile"/><w:rPr>' + $rPr + '</w:rPr></w:pPr>'

# The anchored "KI" logo picture keeps its original (untouched) run formatting.
$drawingRPr = '<w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/>' + `
              '<w:noProof/><w:sz w:val="16"/><w:szCs w:val="16"/>'
$drawingXml = '<w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" ' + `
              'relativeHeight="251661312" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" ' + `
              'wp14:anchorId="331EC226" wp14:editId="7010C1B5"><wp:simplePos x="0" y="0"/>' + `
              '<wp:positionH relativeFrom="column"><wp:posOffset>5655652</wp:posOffset></wp:positionH>' + `
              '<wp:positionV relativeFrom="paragraph"><wp:posOffset>-168812</wp:posOffset></wp:positionV>' + `
              '<wp:extent cx="611945" cy="611945"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/>' + `
              '<wp:docPr id="3" name="Grafik 3" descr="Ein Bild, das Astronomisches Objekt, Kugel, Planet, ' + `
              'Astronomisches Ereignis enth' + [char]0x00E4 + 'lt.&#xA;&#xA;KI-generierte Inhalte k' + `
              [char]0x00F6 + 'nnen fehlerhaft sein."/>' + `
              '<wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' + `
              'noChangeAspect="1"/></wp:cNvGraphicFramePr>' + `
              '<a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">' + `
              '<a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
              '<pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr>' + `
              '<pic:cNvPr id="3" name="Grafik 3" descr="Ein Bild, das Astronomisches Objekt, Kugel, Planet, ' + `
              'Astronomisches Ereignis enth' + [char]0x00E4 + 'lt.&#xA;&#xA;KI-generierte Inhalte k' + `
              [char]0x00F6 + 'nnen fehlerhaft sein."/><pic:cNvPicPr/></pic:nvPicPr>' + `
              '<pic:blipFill><a:blip r:embed="rId1"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}">' + `
              '<a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/>' + `
              '</a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' + `
              '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="611945" cy="611945"/></a:xfrm>' + `
              '<a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData>' + `
              '</a:graphic></wp:anchor></w:drawing>'

function New-Run([string]$text, [bool]$preserveSpace = $false) {
    $space = ""
    if ($preserveSpace) { $space = ' xml:space="preserve"' }
    return '<w:r><w:rPr>' + $rPr + '</w:rPr><w:t' + $space + '>' + $text + '</w:t></w:r>'
}

$apos = [char]0x2019
$ndash = [char]0x2013
$auml = [char]0x00E4

# Paragraph 1: title line, logo, " - KI in der Kita" tagline - all as plain
# runs (no leftover proofErr wrappers).
$p1Open = '<w:p w14:paraId="1470AF61" w14:textId="2D9DFAC2" w:rsidR="0063194B" ' + `
          'w:rsidRPr="00F017E9" w:rsidRDefault="00F017E9" w:rsidP="00F017E9">'
$p1 = $p1Open + $pPr + `
      (New-Run "Don${apos}t Panic Room") + `
      '<w:r><w:rPr>' + $drawingRPr + '</w:rPr>' + $drawingXml + '</w:r>' + `
      (New-Run " $ndash " $true) + `
      (New-Run "KI") + `
      (New-Run " " $true) + `
      (New-Run "in der Kita") + `
      '</w:p>'

# Paragraph 2 (new): "02Station Prompt-Anfaenge Elterngespraeche", with the
# German words wrapped in spell-check proofErr markers.
$p2 = '<w:p>' + $pPr + `
      (New-Run "02S") + `
      (New-Run "tation Prompt-") + `
      '<w:proofErr w:type="spellStart"/>' + (New-Run "Anf${auml}nge") + '<w:proofErr w:type="spellEnd"/>' + `
      (New-Run " " $true) + `
      '<w:proofErr w:type="spellStart"/>' + (New-Run "Elterngespr${auml}che") + '<w:proofErr w:type="spellEnd"/>' + `
      '</w:p>'

$body = $p1 + $p2

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
       '<pkg:part pkg:name="/word/document.xml" ' + `
       'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
       '<pkg:xmlData><w:document ' + `
       'xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" ' + `
       'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' + `
       'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
       'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + `
       'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" ' + `
       'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' + `
       'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
       '<w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$rng.InsertXML($xml)

Write-Host "Header paragraphs now:"
Write-Host $targetHeader.Range.Text
